$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F; the former column F (the "-" header notes column)
# shifts right to become column G, carrying its formatting and contents along.
$ws.Range("F1").EntireColumn.Insert()

# The inserted column F picks up the row formatting from its neighbour for every
# row that had data through column E; strip that back to the default "Normal"
# style for the data rows (the header cell F1 below keeps the inherited style).
$ws.Range("F2:F22").Style = "Normal"

# New header for the inserted column.
$ws.Range("F1").Value = "combo dry g"

# Explain the new column via a note in G2 (a brand-new cell).
$ws.Range("G2").Value = "Combo dry g colum is made up of the sum of the dry g of corresponding number columns"

# Row 2: plain (non-shared) formula summing fine + rhizome dry g.
$ws.Range("F2").Formula = "=SUM(B2, E2)"

# Rows 3-22: fill the same formula down as a shared formula group.
$ws.Range("F3:F22").Formula = "=SUM(B3, E3)"

# Rows where the underlying measurements are notes/missing (or were otherwise
# not carried through the combo column) get "NA" instead of the formula result.
$ws.Range("F8").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("F13").Value = "NA"
$ws.Range("F14").Value = "NA"
$ws.Range("F15:F22").Value = "NA"

# Restore the selection Excel leaves behind after this kind of edit.
[void]$ws.Range("G27").Select()
